$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (BANK OF AMERICA CORP) ---
# Existing Mortgage/Deposit platform-coder-runtime triple (G5:I5) is swapped out
# for a new "Python - API / Vishal / 10" entry, and the previous
# "Python - Selenium / Nimmi" pairing is relocated to the (until now empty)
# M5:O5 triple with a new run time of 60.
$ws.Range("G5").Value = "Python - API"
$ws.Range("H5").Value = "Vishal"
$ws.Range("I5").Value = 10
$ws.Range("M5").Value = "Python - Selenium"
$ws.Range("N5").Value = "Nimmi"
$ws.Range("O5").Value = 60

# --- Row 7 (CITIGROUP INC) ---
# Highlight A7 to match the other highlighted rows (A5/A11/A29): yellow fill.
$ws.Range("A7").Interior.Color = 65535

# --- Row 11 (SYNCHRONY) ---
# Same kind of relocation as row 5: D11:F11 gets a new
# "Python - API / Vishal / 3" entry, and the previous
# "Python - Selenium / Nimmi / 50" triple moves into J11:L11.
$ws.Range("D11").Value = "Python - API"
$ws.Range("E11").Value = "Vishal"
$ws.Range("F11").Value = 3
$ws.Range("J11").Value = "Python - Selenium"
$ws.Range("K11").Value = "Nimmi"
$ws.Range("L11").Value = 50
# Highlight A11 to match the other highlighted rows: yellow fill.
$ws.Range("A11").Interior.Color = 65535

# --- Row 29 (Santander Bank) ---
# Fill in the previously empty Deposit platform-coder-runtime triple.
$ws.Range("D29").Value = "Python - Pagesource"
$ws.Range("E29").Value = "Vishal"
$ws.Range("F29").Value = 10

# --- Selection ---
# The workbook was left with J11 selected.
$null = $ws.Range("J11").Select()
